$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update course load values (columns C, D, E) to reflect the new quarter/course codes
$ws.Range("D3").Value = "2 EF 10A"
$ws.Range("E3").Value = "3 FR 101"

$ws.Range("C5").Value = "3 ER 10A "
$ws.Range("D5").Value = "2 jk 101"
$ws.Range("E5").Value = "1 UK 10A; 1 JKLH 101"

$ws.Range("C7").Value = "2 IC 102"
$ws.Range("D7").Value = "2 BC 102"
$ws.Range("E7").Value = "2 BA 101"

$ws.Range("C8").Value = "2 UP 102"
$ws.Range("D8").Value = "2 PK 101"
$ws.Range("E8").Value = "2 EKP 10A"

$ws.Range("C9").Value = "1 TU 10A, 1 UT 101"
$ws.Range("E9").Value = " 1 EIC 102"

# Update the active selection to D9, matching the new cursor position
$ws.Range("D9").Select()
